$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph "detta kan ske genom npm install och
#    postskripts." and fix the spelling ("postskripts." ->
#    "postscripts", dropping the trailing period), then split the
#    resulting run into three runs: "...och posts" / "c" / "ripts".
# ------------------------------------------------------------------
$target = $d.Content
$ok = $target.Find.Execute(
    "detta kan ske genom npm install och postskripts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "detta kan ske genom npm install och postscripts", 2)

# $target now collapses onto the replaced text (the Find/Replace
# leaves the range spanning the newly inserted text).
$paraRange = $target.Paragraphs(1).Range
$paraStart = $paraRange.Start

# "detta kan ske genom npm install och posts" is 41 characters long.
$splitA = $paraStart + 41
# "c" is the 42nd character.
$splitB = $paraStart + 42

# Toggling a character formatting property on a sub-range forces the
# run to split at the sub-range boundary without altering the visible
# formatting (Bold is set then immediately cleared again).
$subA = $d.Range($paraStart, $splitA)
$subA.Font.Bold = 1
$subA.Font.Bold = 0

$subB = $d.Range($splitA, $splitB)
$subB.Font.Bold = 1
$subB.Font.Bold = 0

# ------------------------------------------------------------------
# 2. Insert the three new bullet paragraphs right after it:
#      ilvl=2  "postscripts"
#      ilvl=3  "installera nödvändiga binaries"
#      ilvl=3  "skapa olika mappar"
#    (all part of the same ListParagraph / numId=6 list as the
#    paragraph above.)
# ------------------------------------------------------------------
$fixedPara = $d.Range($paraStart, $paraStart).Paragraphs(1)
$fixedPara.Range.InsertParagraphAfter()

$p1 = $fixedPara.Next()
$p1.Range.ListFormat.ListLevelNumber = 3
$p1.Range.InsertBefore("postscripts")

$p1 = $fixedPara.Next()
$p1.Range.InsertParagraphAfter()

$p2 = $p1.Next()
$p2.Range.ListFormat.ListLevelNumber = 4
$p2.Range.InsertBefore("installera nödvändiga binaries")

$p2.Range.InsertParagraphAfter()

$p3 = $p2.Next()
$p3.Range.ListFormat.ListLevelNumber = 4
$p3.Range.InsertBefore("skapa olika mappar")

Write-Output "edit applied"
